$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 7.7505571339036123
$ws.Range("C2").Value = 2.7919394720494495
$ws.Range("D2").Value = 1.5104272676029211
$ws.Range("E2").Value = 0.66949841826072998

$ws.Range("B3").Value = 7.0653639578236493
$ws.Range("C3").Value = 12.59800543110876
$ws.Range("D3").Value = 9.7847658980804511
$ws.Range("E3").Value = -3.817925753123717

$ws.Range("B1:E3").Select() | Out-Null
